$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4159.5713
$ws.Range("I43").Value = 3828
$ws.Range("K43").Value = 3828
$ws.Range("M43").Value = -3759

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 2166.2856
$ws.Range("I96").Value = 2024.6
$ws.Range("K96").Value = 6073.799999999999
$ws.Range("M96").Value = -4700.799999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 10563.29
$ws.Range("I132").Value = 7986
$ws.Range("J132").Value = 34618
$ws.Range("K132").Value = 23958
$ws.Range("L132").Value = 103854
$ws.Range("M132").Value = -21428
$ws.Range("N132").Value = -108914

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1836.6129
$ws.Range("I135").Value = 1844.1428
$ws.Range("K135").Value = 16597.2852
$ws.Range("M135").Value = -14062.2852

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 8550.617
$ws.Range("I137").Value = 2986.7646
$ws.Range("K137").Value = 8960.293799999999
$ws.Range("M137").Value = -6410.293799999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2359.11
$ws.Range("J138").Value = 2499.4587
$ws.Range("L138").Value = 7498.3761
$ws.Range("N138").Value = -17778.3761

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2964.9688
$ws.Range("I141").Value = 2539.3
$ws.Range("K141").Value = 7617.900000000001
$ws.Range("M141").Value = -2437.900000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3250.2104
$ws.Range("I32").Value = 2663.3188
$ws.Range("K32").Value = 2663.3188
$ws.Range("M32").Value = -2376.3188

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 19980.285
$ws.Range("I99").Value = 21029.904
$ws.Range("J99").Value = 16831.428
$ws.Range("K99").Value = 21029.904
$ws.Range("L99").Value = 16831.428
$ws.Range("M99").Value = -19531.904
$ws.Range("N99").Value = -19827.428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 22964.121
$ws.Range("I134").Value = 23226.617
$ws.Range("K134").Value = 69679.851
$ws.Range("M134").Value = -67144.851

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8165.2144
$ws.Range("I16").Value = 5831.5
$ws.Range("K16").Value = 5831.5
$ws.Range("M16").Value = -5544.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 905.92
$ws.Range("I22").Value = 551.2353000000001
$ws.Range("K22").Value = 551.2353000000001
$ws.Range("M22").Value = -201.2353000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24454.902
$ws.Range("I31").Value = 17620
$ws.Range("K31").Value = 17620
$ws.Range("M31").Value = -17325

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 24454.902
$ws.Range("I34").Value = 17620
$ws.Range("K34").Value = 17620
$ws.Range("M34").Value = -17418

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 7405.3125
$ws.Range("I105").Value = 9207.166999999999
$ws.Range("K105").Value = 9207.166999999999
$ws.Range("M105").Value = -7460.166999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 8165.2144
$ws.Range("I113").Value = 5831.5
$ws.Range("K113").Value = 5831.5
$ws.Range("M113").Value = -3661.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2375
$ws.Range("I122").Value = 2375
$ws.Range("K122").Value = 7125
$ws.Range("M122").Value = -4675

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7489.278
$ws.Range("I132").Value = 1674.5625
$ws.Range("K132").Value = 5023.6875
$ws.Range("M132").Value = -2493.6875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1478.75
$ws.Range("I26").Value = 2897.5
$ws.Range("J26").Value = 60
$ws.Range("K26").Value = 8692.5
$ws.Range("L26").Value = 180
$ws.Range("M26").Value = -8404.5
$ws.Range("N26").Value = -756

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 2033.1666
$ws.Range("I60").Value = 2239.8
$ws.Range("J60").Value = 1000
$ws.Range("K60").Value = 6719.400000000001
$ws.Range("L60").Value = 3000
$ws.Range("M60").Value = -6468.400000000001
$ws.Range("N60").Value = -3502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 671
$ws.Range("J68").Value = 625.25
$ws.Range("L68").Value = 1875.75
$ws.Range("N68").Value = -3497.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 671
$ws.Range("J71").Value = 625.25
$ws.Range("L71").Value = 5627.25
$ws.Range("N71").Value = -13739.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 566.44446
$ws.Range("I92").Value = 580
$ws.Range("J92").Value = 549.5
$ws.Range("K92").Value = 1740
$ws.Range("L92").Value = 1648.5
$ws.Range("M92").Value = -492
$ws.Range("N92").Value = -4144.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 16973.25
$ws.Range("J98").Value = 17949
$ws.Range("L98").Value = 53847
$ws.Range("N98").Value = -56843

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 3705606.2
$ws.Range("I109").Value = 2140.625
$ws.Range("K109").Value = 6421.875
$ws.Range("M109").Value = -5381.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2567.5908
$ws.Range("I102").Value = 2892
$ws.Range("K102").Value = 2892
$ws.Range("M102").Value = -1270

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 800.1667
$ws.Range("I107").Value = 546.875
$ws.Range("J107").Value = 1306.75
$ws.Range("K107").Value = 546.875
$ws.Range("L107").Value = 1306.75
$ws.Range("M107").Value = 1373.125
$ws.Range("N107").Value = -5146.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 13177.29
$ws.Range("I132").Value = 10936.091
$ws.Range("K132").Value = 32808.273
$ws.Range("M132").Value = -30278.273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5549.3794
$ws.Range("I7").Value = 6128.727
$ws.Range("K7").Value = 6128.727
$ws.Range("M7").Value = -6016.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1705.2778
$ws.Range("I16").Value = 1761.4706
$ws.Range("J16").Value = 750
$ws.Range("K16").Value = 1761.4706
$ws.Range("L16").Value = 750
$ws.Range("M16").Value = -1591.4706
$ws.Range("N16").Value = -1090

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3681.5293
$ws.Range("I82").Value = 4309.778
$ws.Range("J82").Value = 2974.75
$ws.Range("K82").Value = 4309.778
$ws.Range("L82").Value = 2974.75
$ws.Range("M82").Value = -3948.778
$ws.Range("N82").Value = -3696.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3681.5293
$ws.Range("I85").Value = 4309.778
$ws.Range("J85").Value = 2974.75
$ws.Range("K85").Value = 4309.778
$ws.Range("L85").Value = 2974.75
$ws.Range("M85").Value = -3061.778
$ws.Range("N85").Value = -5470.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3330.9736
$ws.Range("I93").Value = 1620.1111
$ws.Range("K93").Value = 1620.1111
$ws.Range("M93").Value = -372.1111000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2564.0417
$ws.Range("I100").Value = 2257.6924
$ws.Range("J100").Value = 2926.0908
$ws.Range("K100").Value = 2257.6924
$ws.Range("L100").Value = 2926.0908
$ws.Range("M100").Value = -1716.6924
$ws.Range("N100").Value = -4008.0908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5549.3794
$ws.Range("I126").Value = 6128.727
$ws.Range("K126").Value = 18386.181
$ws.Range("M126").Value = -15916.181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1391.1666
$ws.Range("J96").Value = 1490.4
$ws.Range("L96").Value = 1490.4
$ws.Range("N96").Value = -4236.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 659.7059
$ws.Range("I100").Value = 568.25
$ws.Range("K100").Value = 1136.5
$ws.Range("M100").Value = -595.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 22666.666
$ws.Range("J109").Value = 22666.666
$ws.Range("L109").Value = 22666.666
$ws.Range("N109").Value = -25440.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7462.8335
$ws.Range("I132").Value = 2544.6
$ws.Range("K132").Value = 7633.799999999999
$ws.Range("M132").Value = -5103.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 9450.880999999999
$ws.Range("I136").Value = 1095.3226
$ws.Range("J136").Value = 32998.363
$ws.Range("K136").Value = 3285.9678
$ws.Range("L136").Value = 98995.08899999999
$ws.Range("M136").Value = -735.9677999999999
$ws.Range("N136").Value = -104095.089

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 137081
$ws.Range("J140").Value = 137081
$ws.Range("L140").Value = 137081
$ws.Range("N140").Value = -147441
